# Weekly update: a new week's price record is inserted at row 34,
# pushing all subsequent records down by one row (old row 156 becomes
# the new row 157).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 34; Excel shifts rows 34:156 down to 35:157
# and extends the used range to row 157 automatically.
$ws.Rows.Item(34).Insert()

# Populate the constant columns (identical on every data row in this sheet)
# by copying them from the row directly below (the row that used to be 34).
$ws.Range("A34").Value = $ws.Range("A35").Value()
$ws.Range("B34").Value = $ws.Range("B35").Value()
$ws.Range("C34").Value = $ws.Range("C35").Value()
$ws.Range("E34").Value = $ws.Range("E35").Value()
$ws.Range("F34").Value = $ws.Range("F35").Value()
$ws.Range("G34").Value = $ws.Range("G35").Value()
$ws.Range("H34").Value = $ws.Range("H35").Value()
$ws.Range("N34").Value = $ws.Range("N35").Value()
$ws.Range("Q34").Value = $ws.Range("Q35").Value()
$ws.Range("R34").Value = $ws.Range("R35").Value()

# New record's own data.
$ws.Range("D34").Value = 44453
$ws.Range("I34").Value = "Primera"
$ws.Range("J34").Value = 160
$ws.Range("K34").Value = 700
$ws.Range("L34").Value = 750
$ws.Range("M34").Value = 725
$ws.Range("O34").Value = "Región del Maule"
$ws.Range("P34").Value = 725
